$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1379.8
$ws.Cells.Item(32, 9).Value = 1349.5
$ws.Cells.Item(32, 10).Value = 1400
$ws.Cells.Item(32, 11).Value = 1349.5
$ws.Cells.Item(32, 12).Value = 1400
$ws.Cells.Item(32, 13).Value = -1023.5
$ws.Cells.Item(32, 14).Value = -2052
$ws.Cells.Item(33, 8).Value = 452.6316
$ws.Cells.Item(33, 9).Value = 372.33334
$ws.Cells.Item(33, 10).Value = 753.75
$ws.Cells.Item(33, 11).Value = 372.33334
$ws.Cells.Item(33, 12).Value = 753.75
$ws.Cells.Item(33, 13).Value = -143.33334
$ws.Cells.Item(33, 14).Value = -1211.75
$ws.Cells.Item(80, 8).Value = 2963774.5
$ws.Cells.Item(80, 9).Value = 1358.3334
$ws.Cells.Item(80, 10).Value = 3852499.2
$ws.Cells.Item(80, 11).Value = 4075.0002
$ws.Cells.Item(80, 12).Value = 11557497.6
$ws.Cells.Item(80, 13).Value = -3077.0002
$ws.Cells.Item(80, 14).Value = -11559493.6
$ws.Cells.Item(83, 8).Value = 2963774.5
$ws.Cells.Item(83, 9).Value = 1358.3334
$ws.Cells.Item(83, 10).Value = 3852499.2
$ws.Cells.Item(83, 11).Value = 12225.0006
$ws.Cells.Item(83, 12).Value = 34672492.8
$ws.Cells.Item(83, 13).Value = -7233.000599999999
$ws.Cells.Item(83, 14).Value = -34682476.8
$ws.Cells.Item(107, 8).Value = 461.53333
$ws.Cells.Item(107, 9).Value = 280
$ws.Cells.Item(107, 10).Value = 552.3
$ws.Cells.Item(107, 11).Value = 280
$ws.Cells.Item(107, 12).Value = 552.3
$ws.Cells.Item(107, 13).Value = 1640
$ws.Cells.Item(107, 14).Value = -4392.3
$ws.Cells.Item(112, 8).Value = 3907296.2
$ws.Cells.Item(112, 10).Value = 1099.3
$ws.Cells.Item(112, 12).Value = 3297.9
$ws.Cells.Item(112, 14).Value = -5513.9
$ws.Cells.Item(129, 8).Value = 899.13336
$ws.Cells.Item(129, 10).Value = 911.2241
$ws.Cells.Item(129, 12).Value = 2733.6723
$ws.Cells.Item(129, 14).Value = -12733.6723
$ws.Cells.Item(138, 8).Value = 2770.025
$ws.Cells.Item(138, 9).Value = 1767.3125
$ws.Cells.Item(138, 10).Value = 3438.5
$ws.Cells.Item(138, 11).Value = 5301.9375
$ws.Cells.Item(138, 12).Value = 10315.5
$ws.Cells.Item(138, 13).Value = -161.9375
$ws.Cells.Item(138, 14).Value = -20595.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 26248.715
$ws.Cells.Item(132, 9).Value = 2013.6666
$ws.Cells.Item(132, 11).Value = 6040.9998
$ws.Cells.Item(132, 13).Value = -3510.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1131.3334
$ws.Cells.Item(20, 9).Value = 997
$ws.Cells.Item(20, 10).Value = 1400
$ws.Cells.Item(20, 11).Value = 997
$ws.Cells.Item(20, 12).Value = 1400
$ws.Cells.Item(20, 13).Value = -750
$ws.Cells.Item(20, 14).Value = -1894
$ws.Cells.Item(134, 8).Value = 3084.4333
$ws.Cells.Item(134, 9).Value = 3501.4167
$ws.Cells.Item(134, 10).Value = 1416.5
$ws.Cells.Item(134, 11).Value = 10504.2501
$ws.Cells.Item(134, 12).Value = 4249.5
$ws.Cells.Item(134, 13).Value = -7969.250100000001
$ws.Cells.Item(134, 14).Value = -9319.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1128.2142
$ws.Cells.Item(16, 9).Value = 1049.5834
$ws.Cells.Item(16, 11).Value = 1049.5834
$ws.Cells.Item(16, 13).Value = -762.5834
$ws.Cells.Item(94, 8).Value = 3088.158
$ws.Cells.Item(94, 9).Value = 2085.7144
$ws.Cells.Item(94, 10).Value = 3672.9167
$ws.Cells.Item(94, 11).Value = 2085.7144
$ws.Cells.Item(94, 12).Value = 3672.9167
$ws.Cells.Item(94, 13).Value = -1634.7144
$ws.Cells.Item(94, 14).Value = -4574.9167
$ws.Cells.Item(113, 8).Value = 1128.2142
$ws.Cells.Item(113, 9).Value = 1049.5834
$ws.Cells.Item(113, 11).Value = 1049.5834
$ws.Cells.Item(113, 13).Value = 1120.4166
$ws.Cells.Item(132, 8).Value = 3729.7058
$ws.Cells.Item(132, 9).Value = 1236.25
$ws.Cells.Item(132, 11).Value = 3708.75
$ws.Cells.Item(132, 13).Value = -1178.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 1700
$ws.Cells.Item(69, 10).Value = 1800
$ws.Cells.Item(69, 12).Value = 5400
$ws.Cells.Item(69, 14).Value = -7022
$ws.Cells.Item(72, 8).Value = 1700
$ws.Cells.Item(72, 10).Value = 1800
$ws.Cells.Item(72, 12).Value = 16200
$ws.Cells.Item(72, 14).Value = -24312
$ws.Cells.Item(80, 8).Value = 2324.875
$ws.Cells.Item(80, 9).Value = 1999.6666
$ws.Cells.Item(80, 10).Value = 2520
$ws.Cells.Item(80, 11).Value = 5998.9998
$ws.Cells.Item(80, 12).Value = 7560
$ws.Cells.Item(80, 13).Value = -5062.9998
$ws.Cells.Item(80, 14).Value = -9432
$ws.Cells.Item(83, 8).Value = 2324.875
$ws.Cells.Item(83, 9).Value = 1999.6666
$ws.Cells.Item(83, 10).Value = 2520
$ws.Cells.Item(83, 11).Value = 17996.9994
$ws.Cells.Item(83, 12).Value = 22680
$ws.Cells.Item(83, 13).Value = -13316.9994
$ws.Cells.Item(83, 14).Value = -32040
$ws.Cells.Item(131, 8).Value = 801.13
$ws.Cells.Item(131, 10).Value = 824.97894
$ws.Cells.Item(131, 12).Value = 2474.93682
$ws.Cells.Item(131, 14).Value = -12554.93682

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1823.1428
$ws.Cells.Item(102, 9).Value = 1768.6666
$ws.Cells.Item(102, 10).Value = 2150
$ws.Cells.Item(102, 11).Value = 1768.6666
$ws.Cells.Item(102, 12).Value = 2150
$ws.Cells.Item(102, 13).Value = -146.6666
$ws.Cells.Item(102, 14).Value = -5394
$ws.Cells.Item(132, 8).Value = 56119.8
$ws.Cells.Item(132, 9).Value = 7200
$ws.Cells.Item(132, 10).Value = 129499.5
$ws.Cells.Item(132, 11).Value = 21600
$ws.Cells.Item(132, 12).Value = 388498.5
$ws.Cells.Item(132, 13).Value = -19070
$ws.Cells.Item(132, 14).Value = -393558.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1059
$ws.Cells.Item(22, 9).Value = 797.5
$ws.Cells.Item(22, 10).Value = 1233.3334
$ws.Cells.Item(22, 11).Value = 797.5
$ws.Cells.Item(22, 12).Value = 1233.3334
$ws.Cells.Item(22, 13).Value = -502.5
$ws.Cells.Item(22, 14).Value = -1823.3334
$ws.Cells.Item(27, 8).Value = 1059
$ws.Cells.Item(27, 9).Value = 797.5
$ws.Cells.Item(27, 10).Value = 1233.3334
$ws.Cells.Item(27, 11).Value = 797.5
$ws.Cells.Item(27, 12).Value = 1233.3334
$ws.Cells.Item(27, 13).Value = -690.5
$ws.Cells.Item(27, 14).Value = -1447.3334
$ws.Cells.Item(46, 8).Value = 777.2174
$ws.Cells.Item(46, 9).Value = 694.0952
$ws.Cells.Item(46, 11).Value = 694.0952
$ws.Cells.Item(46, 13).Value = -506.0952
$ws.Cells.Item(55, 8).Value = 260.5263
$ws.Cells.Item(55, 9).Value = 190
$ws.Cells.Item(55, 10).Value = 279.33334
$ws.Cells.Item(55, 11).Value = 190
$ws.Cells.Item(55, 12).Value = 279.33334
$ws.Cells.Item(55, 13).Value = -17
$ws.Cells.Item(55, 14).Value = -625.33334

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2446
$ws.Cells.Item(132, 9).Value = 2063.8
$ws.Cells.Item(132, 10).Value = 3083
$ws.Cells.Item(132, 11).Value = 6191.400000000001
$ws.Cells.Item(132, 12).Value = 9249
$ws.Cells.Item(132, 13).Value = -3661.400000000001
$ws.Cells.Item(132, 14).Value = -14309

